$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.39011287689209
$ws.Range("B1").Value = 2.686439514160156
$ws.Range("C1").Value = 2.042529821395874
$ws.Range("D1").Value = 1.908952355384827
$ws.Range("E1").Value = 1.943543910980225
